# Updated cryptos list on Fri Sep 15 09:24:04 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) columns with latest scraped values,
# and corrects the swapped Polkadot / WrappedEther rows (13 & 14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values that look like plain numbers must be forced to Text so
# Excel doesn't silently reformat them (e.g. "1.50" -> 1.5, "0.0512" -> 5.12E-2),
# matching the original inline-string cell content byte-for-byte.

$ws.Range("D2").Value = "26.661.81"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.630.95"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.32"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.22"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("D12").Value = "1.858.30"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.604.49"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "26.645.59"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.47"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.50"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.36"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.71"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "1.214.69"
$ws.Range("E36").Value = "  +4.85%  "
$ws.Range("E37").Value = "  +4.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.804"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.35"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "1.767.62"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.74"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.58"
$ws.Range("E50").Value = "  +3.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.410"
$ws.Range("E51").Value = "  -0.04%  "
